# Insert a new data row at row 135 (pushing existing rows 135..212 down to 136..213)
# and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(135).Insert()

$ws.Cells.Item(135, 1).Value = 1
$ws.Cells.Item(135, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(135, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(135, 4).Value = 44529
$ws.Cells.Item(135, 5).Value = 15
$ws.Cells.Item(135, 6).Value = 100114013
$ws.Cells.Item(135, 7).Value = "Zanahoria"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 90
$ws.Cells.Item(135, 11).Value = 15000
$ws.Cells.Item(135, 12).Value = 16000
$ws.Cells.Item(135, 13).Value = 15500
$ws.Cells.Item(135, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(135, 15).Value = "Valle de Camiña"
$ws.Cells.Item(135, 16).Value = 620
$ws.Cells.Item(135, 17).Value = 25
$ws.Cells.Item(135, 18).Value = "Hortaliza"
